$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 171.5
$ws.Range("I12").Value = 185.6
$ws.Range("K12").Value = 185.6
$ws.Range("M12").Value = -15.59999999999999
$ws.Range("H17").Value = 912.53845
$ws.Range("J17").Value = 934.55554
$ws.Range("L17").Value = 2803.66662
$ws.Range("N17").Value = -3139.66662
$ws.Range("H18").Value = 237.42857
$ws.Range("I18").Value = 254
$ws.Range("J18").Value = 196
$ws.Range("K18").Value = 254
$ws.Range("L18").Value = 196
$ws.Range("M18").Value = 30
$ws.Range("N18").Value = -764
$ws.Range("H32").Value = 6737.52
$ws.Range("J32").Value = 6016
$ws.Range("L32").Value = 6016
$ws.Range("N32").Value = -6668
$ws.Range("H43").Value = 5290.5454
$ws.Range("I43").Value = 4264.6665
$ws.Range("J43").Value = 5675.25
$ws.Range("K43").Value = 4264.6665
$ws.Range("L43").Value = 5675.25
$ws.Range("M43").Value = -4195.6665
$ws.Range("N43").Value = -5813.25
$ws.Range("H62").Value = 10013.625
$ws.Range("I62").Value = 8544.416999999999
$ws.Range("K62").Value = 8544.416999999999
$ws.Range("M62").Value = -7920.416999999999
$ws.Range("H65").Value = 10013.625
$ws.Range("I65").Value = 8544.416999999999
$ws.Range("K65").Value = 42722.085
$ws.Range("M65").Value = -39602.085
$ws.Range("H111").Value = 200
$ws.Range("I111").Value = 200
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 600
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("M111").Value = 2467
$ws.Range("H112").Value = 13440.6
$ws.Range("I112").Value = 650
$ws.Range("J112").Value = 16638.25
$ws.Range("K112").Value = 1950
$ws.Range("L112").Value = 49914.75
$ws.Range("M112").Value = -842
$ws.Range("N112").Value = -52130.75
$ws.Range("H116").Value = 16704.5
$ws.Range("I116").Value = 17906.334
$ws.Range("J116").Value = 11897.167
$ws.Range("K116").Value = 17906.334
$ws.Range("L116").Value = 11897.167
$ws.Range("M116").Value = -14464.334
$ws.Range("N116").Value = -18781.167
$ws.Range("H135").Value = 5002.5
$ws.Range("I135").Value = 4671.6665
$ws.Range("K135").Value = 42044.9985
$ws.Range("M135").Value = -39509.9985
$ws.Range("H138").Value = 69358.13
$ws.Range("I138").Value = 2883.7856
$ws.Range("J138").Value = 999999
$ws.Range("K138").Value = 8651.356800000001
$ws.Range("L138").Value = 2999997
$ws.Range("M138").Value = -3511.356800000001
$ws.Range("N138").Value = -3010277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 673.1667
$ws.Range("I2").Value = 633.9091
$ws.Range("K2").Value = 633.9091
$ws.Range("M2").Value = -520.9091
$ws.Range("H61").Value = 10089.583
$ws.Range("I61").Value = 1345.25
$ws.Range("K61").Value = 1345.25
$ws.Range("M61").Value = -1133.25
$ws.Range("H74").Value = 383490.06
$ws.Range("I74").Value = 601360.1
$ws.Range("K74").Value = 601360.1
$ws.Range("M74").Value = -600486.1
$ws.Range("H77").Value = 383490.06
$ws.Range("I77").Value = 601360.1
$ws.Range("K77").Value = 3006800.5
$ws.Range("M77").Value = -3002432.5
$ws.Range("H116").Value = 673.1667
$ws.Range("I116").Value = 633.9091
$ws.Range("K116").Value = 633.9091
$ws.Range("M116").Value = 1660.0909
$ws.Range("H122").Value = 2940.5
$ws.Range("J122").Value = 3411.111
$ws.Range("L122").Value = 10233.333
$ws.Range("N122").Value = -15133.333
$ws.Range("H136").Value = 10089.583
$ws.Range("I136").Value = 1345.25
$ws.Range("K136").Value = 4035.75
$ws.Range("M136").Value = -1485.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 673.1667
$ws.Range("I3").Value = 633.9091
$ws.Range("K3").Value = 633.9091
$ws.Range("M3").Value = -519.9091
$ws.Range("H64").Value = 1729
$ws.Range("I64").Value = 1756
$ws.Range("J64").Value = 1715.5
$ws.Range("K64").Value = 1756
$ws.Range("L64").Value = 1715.5
$ws.Range("M64").Value = -1531
$ws.Range("N64").Value = -2165.5
$ws.Range("H67").Value = 1729
$ws.Range("I67").Value = 1756
$ws.Range("J67").Value = 1715.5
$ws.Range("K67").Value = 1756
$ws.Range("L67").Value = 1715.5
$ws.Range("M67").Value = -976
$ws.Range("N67").Value = -3275.5
$ws.Range("H86").Value = 2199.6667
$ws.Range("I86").Value = 1699.5714
$ws.Range("J86").Value = 2899.8
$ws.Range("K86").Value = 1699.5714
$ws.Range("L86").Value = 2899.8
$ws.Range("M86").Value = -576.5714
$ws.Range("N86").Value = -5145.8
$ws.Range("H89").Value = 2199.6667
$ws.Range("I89").Value = 1699.5714
$ws.Range("J89").Value = 2899.8
$ws.Range("K89").Value = 8497.857
$ws.Range("L89").Value = 14499
$ws.Range("M89").Value = -2881.857
$ws.Range("N89").Value = -25731

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 5002.3335
$ws.Range("I94").Value = 4498
$ws.Range("J94").Value = 5254.5
$ws.Range("K94").Value = 4498
$ws.Range("L94").Value = 5254.5
$ws.Range("M94").Value = -4047
$ws.Range("N94").Value = -6156.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 176916560
$ws.Range("J4").Value = 239861420
$ws.Range("L4").Value = 719584260
$ws.Range("N4").Value = -719584484
$ws.Range("H50").Value = 25350
$ws.Range("I50").Value = 25350
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 76050
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("M50").Value = -75569
$ws.Range("H53").Value = 25350
$ws.Range("I53").Value = 25350
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 76050
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("M53").Value = -75569
$ws.Range("H63").Value = 17903.666
$ws.Range("I63").Value = 17903.666
$ws.Range("K63").Value = 53710.99800000001
$ws.Range("M63").Value = -52961.99800000001
$ws.Range("H64").Value = 8548.5
$ws.Range("J64").Value = 4310.7144
$ws.Range("L64").Value = 12932.1432
$ws.Range("N64").Value = -13472.1432
$ws.Range("H66").Value = 17903.666
$ws.Range("I66").Value = 17903.666
$ws.Range("K66").Value = 161132.994
$ws.Range("M66").Value = -157388.994
$ws.Range("H67").Value = 8548.5
$ws.Range("J67").Value = 4310.7144
$ws.Range("L67").Value = 12932.1432
$ws.Range("N67").Value = -14804.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 603666.5
$ws.Range("I3").Value = 278250
$ws.Range("J3").Value = 1254499.5
$ws.Range("K3").Value = 278250
$ws.Range("L3").Value = 1254499.5
$ws.Range("M3").Value = -278134
$ws.Range("N3").Value = -1254731.5
$ws.Range("H9").Value = 1333.3334
$ws.Range("I9").Value = 1333.3334
$ws.Range("K9").Value = 1333.3334
$ws.Range("M9").Value = -1163.3334
$ws.Range("H10").Value = 2777.3333
$ws.Range("I10").Value = 549.5
$ws.Range("J10").Value = 3891.25
$ws.Range("K10").Value = 549.5
$ws.Range("L10").Value = 3891.25
$ws.Range("M10").Value = -380.5
$ws.Range("N10").Value = -4229.25
$ws.Range("H11").Value = 167609.08
$ws.Range("I11").Value = 222300
$ws.Range("J11").Value = 3536.3333
$ws.Range("K11").Value = 222300
$ws.Range("L11").Value = 3536.3333
$ws.Range("M11").Value = -222161
$ws.Range("N11").Value = -3814.3333
$ws.Range("H12").Value = 10000
$ws.Range("I12").Value = 10000
$ws.Range("K12").Value = 10000
$ws.Range("M12").Value = -9860
$ws.Range("H14").Value = 582020160
$ws.Range("I14").Value = 667000000
$ws.Range("J14").Value = 454550430
$ws.Range("K14").Value = 667000000
$ws.Range("L14").Value = 454550430
$ws.Range("M14").Value = -666999832
$ws.Range("N14").Value = -454550766
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H29").Value = 15288.8
$ws.Range("I29").Value = 5000
$ws.Range("J29").Value = 20433.2
$ws.Range("K29").Value = 5000
$ws.Range("L29").Value = 20433.2
$ws.Range("M29").Value = -4710
$ws.Range("N29").Value = -21013.2
$ws.Range("H102").Value = 2668.0715
$ws.Range("I102").Value = 2279.4167
$ws.Range("K102").Value = 2279.4167
$ws.Range("M102").Value = -657.4167000000002
$ws.Range("H122").Value = 21741776
$ws.Range("J122").Value = 55558756
$ws.Range("L122").Value = 166676268
$ws.Range("N122").Value = -166681168
$ws.Range("H132").Value = 3001.7222
$ws.Range("I132").Value = 2675.6
$ws.Range("K132").Value = 8026.799999999999
$ws.Range("M132").Value = -5496.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1825.25
$ws.Range("I40").Value = 1825.25
$ws.Range("K40").Value = 1825.25
$ws.Range("M40").Value = -1689.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 723.5833
$ws.Range("I107").Value = 609.05554
$ws.Range("J107").Value = 1067.1666
$ws.Range("K107").Value = 1827.16662
$ws.Range("L107").Value = 3201.4998
$ws.Range("M107").Value = 92.83338000000003
$ws.Range("N107").Value = -7041.4998
